# Apply cryptos list price/volume/coin-row updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing text storage (avoids Excel auto-
# converting numeric-looking strings like "0.5228" into real numbers), then
# restores the cell to its original (default/"Normal") style so no stray
# number-format style is left behind on the cell.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "26.091.05"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "1.650.51"
$ws.Range("E3").Value = "  -0.97%  "
Set-TextValue "D4" "1.006"
$ws.Range("E4").Value = "  -0.38%  "
Set-TextValue "D5" "218.22"
$ws.Range("E5").Value = "  -0.18%  "
Set-TextValue "D6" "0.5228"
$ws.Range("E6").Value = "  -1.82%  "
Set-TextValue "D7" "1.006"
$ws.Range("E7").Value = "  -0.40%  "
Set-TextValue "D8" "0.2612"
$ws.Range("E8").Value = "  -1.02%  "
Set-TextValue "D9" "0.06267"
$ws.Range("E9").Value = "  -1.81%  "
Set-TextValue "D10" "20.45"
$ws.Range("E10").Value = "  -0.48%  "
Set-TextValue "D11" "0.07816"
$ws.Range("E11").Value = "  +0.00%  "
Set-TextValue "D12" "4.464"
$ws.Range("E12").Value = "  -2.11%  "
$ws.Range("D13").Value = "1.666.45"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "1.884.33"
$ws.Range("E14").Value = "  -0.59%  "
Set-TextValue "D15" "0.5517"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").Value = "0.0₅7967"
$ws.Range("E16").Value = "  -2.79%  "
Set-TextValue "D17" "64.86"
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("D18").Value = "26.120.86"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("E19").Value = "  -0.51%  "
Set-TextValue "D20" "4.615"
$ws.Range("E20").Value = "  -1.40%  "
Set-TextValue "D21" "194.36"
$ws.Range("E21").Value = "  +0.47%  "
Set-TextValue "D22" "10.03"
$ws.Range("E22").Value = "  -1.80%  "
Set-TextValue "D23" "5.936"
$ws.Range("E23").Value = "  -1.61%  "
Set-TextValue "D24" "1.008"
$ws.Range("E24").Value = "  -0.33%  "
Set-TextValue "D25" "146.44"
$ws.Range("E25").Value = "  +0.53%  "
Set-TextValue "D26" "0.1201"
$ws.Range("E26").Value = "  -2.10%  "
Set-TextValue "D27" "7.142"
$ws.Range("E27").Value = "  -0.76%  "
Set-TextValue "D28" "15.91"
$ws.Range("E28").Value = "  -1.24%  "
Set-TextValue "D29" "1.478"
$ws.Range("E29").Value = "  -0.31%  "
Set-TextValue "D30" "0.05684"
$ws.Range("E30").Value = "  -3.50%  "
Set-TextValue "D31" "1.270"
$ws.Range("E31").Value = "  -0.74%  "
Set-TextValue "D32" "3.467"
$ws.Range("E32").Value = "  -3.71%  "
Set-TextValue "D33" "3.319"
$ws.Range("E33").Value = "  +1.30%  "
Set-TextValue "D34" "1.578"
$ws.Range("E34").Value = "  -2.22%  "
Set-TextValue "D35" "2.794"
$ws.Range("E35").Value = "  -1.15%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D36" "2.416"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D37" "0.9461"
$ws.Range("E37").Value = "  -1.81%  "
Set-TextValue "D38" "0.5658"
$ws.Range("E38").Value = "  -2.34%  "
Set-TextValue "D39" "0.01590"
$ws.Range("E39").Value = "  -0.98%  "
Set-TextValue "D40" "5.921"
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("D41").Value = "1.060.55"
$ws.Range("E41").Value = "  +1.20%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D42" "1.006"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D43" "0.8410"
$ws.Range("E43").Value = "  -2.68%  "
Set-TextValue "D44" "103.07"
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("D45").Value = "1.795.42"
$ws.Range("E45").Value = "  -0.57%  "
Set-TextValue "D46" "57.35"
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈107"
$ws.Range("E47").Value = "  +1.64%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D48" "0.05391"
$ws.Range("E48").Value = "  +4.47%  "
$ws.Range("B49").Value = "Frax"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D49" "1.006"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D50" "0.4398"
$ws.Range("E50").Value = "  +0.42%  "
Set-TextValue "D51" "7.990"
$ws.Range("E51").Value = "  -0.59%  "

Write-Host "Applied 115 cell updates"
